# Adds 10 new shipment/log rows (rows 5-14) to the "lojistik" sheet,
# extending the used range from A1:M4 to A1:M14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then the 13 column values (A..M).
# Columns A,B,C,E,G,I,K are text; D,F,H,J,L,M are numeric.
$rows = @(
    @{ R = 5;  A = "31.03.2026"; B = "10:00"; C = "55NM123"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Nisa Karaman";  L = 9;   M = 10 },
    @{ R = 6;  A = "31.03.2026"; B = "10:00"; C = "55NM123"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Nisa Karaman";  L = 9;   M = 10 },
    @{ R = 7;  A = "09.10.1998"; B = "10:00"; C = "55NM123"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Nisa Karaman";  L = 9;   M = 10 },
    @{ R = 8;  A = "00.10.1998"; B = "10:00"; C = "55NM123"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Nisa Karaman";  L = 9;   M = 10 },
    @{ R = 9;  A = "00.10.42";   B = "10:00"; C = "55NM123"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Nisa Karaman";  L = 9;   M = 10 },
    @{ R = 10; A = "03.05.1979"; B = "10:00"; C = "55NM123"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Nisa Karaman";  L = 9;   M = 10 },
    @{ R = 11; A = "03.05.1979"; B = "10:00"; C = "55NM123"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Nisa Karaman";  L = 9;   M = 10 },
    @{ R = 12; A = "123";        B = "2357";  C = "2578";    D = 853; E = "2368";  F = 4680; G = "3568";  H = 169; I = "9643";  J = 3827; K = "Sgkhfbll";      L = 588; M = 689 },
    @{ R = 13; A = "15.02.2025"; B = "09:00"; C = "68HS574"; D = 5;   E = "11:30"; F = 6;    G = "13:45"; H = 14;  I = "14:30"; J = 1;    K = "Melih Karaman"; L = 9;   M = 10 },
    @{ R = 14; A = "11.07.2025"; B = "10:00"; C = "45HD132"; D = 564; E = "12:00"; F = 614;  G = "12:20"; H = 617; I = "14:30"; J = 50;   K = "Ela karaman ";  L = 116; M = 23 }
)

# Columns that hold text even when the value looks purely numeric
# (dates, times, plate numbers, names).
$textCols = @("A", "B", "C", "E", "G", "I", "K")
$numCols  = @("D", "F", "H", "J", "L", "M")

foreach ($row in $rows) {
    $r = $row.R

    foreach ($col in $textCols) {
        $value = [string]$row[$col]
        $cell = $ws.Range("$col$r")

        # Excel auto-coerces cell input that "looks like" a number or a
        # date into that type. Our source values must stay literal text
        # (they are free-form log entries, not real numbers/dates), so
        # pre-format any value that would otherwise be auto-converted.
        $looksNumeric = $value -match '^[0-9]+$'
        $looksDate = $false
        if ($value -match '^(\d{1,2})\.(\d{1,2})\.(\d{2,4})$') {
            $d = [int]$Matches[1]
            $m = [int]$Matches[2]
            if ($d -ge 1 -and $d -le 12 -and $m -ge 1 -and $m -le 12) {
                $looksDate = $true
            }
        }

        if ($looksNumeric -or $looksDate) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
    }

    foreach ($col in $numCols) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}

Write-Output "Added rows 5-14 to sheet"
